$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number-format on the D/E columns we touch so that Excel
# keeps numeric-looking values (e.g. "226.26") as literal text,
# matching the workbook's original inlineStr cells instead of
# letting COM auto-coerce them into floating point numbers.
$cells = @()
$cells += "D2"
$cells += "E2"
$cells += "D3"
$cells += "E3"
$cells += "E4"
$cells += "D5"
$cells += "E5"
$cells += "E6"
$cells += "E7"
$cells += "D8"
$cells += "E8"
$cells += "E9"
$cells += "D10"
$cells += "E10"
$cells += "E11"
$cells += "E12"
$cells += "D13"
$cells += "E13"
$cells += "D14"
$cells += "E14"
$cells += "E15"
$cells += "D16"
$cells += "E16"
$cells += "E17"
$cells += "E18"
$cells += "D19"
$cells += "E19"
$cells += "E20"
$cells += "D21"
$cells += "E21"
$cells += "E22"
$cells += "D23"
$cells += "E23"
$cells += "D24"
$cells += "E24"
$cells += "D25"
$cells += "E26"
$cells += "D27"
$cells += "E27"
$cells += "E28"
$cells += "E29"
$cells += "E30"
$cells += "D31"
$cells += "E31"
$cells += "D32"
$cells += "E32"
$cells += "E33"
$cells += "E34"
$cells += "D35"
$cells += "E35"
$cells += "D36"
$cells += "E36"
$cells += "D37"
$cells += "E37"
$cells += "D38"
$cells += "E38"
$cells += "E39"
$cells += "E40"
$cells += "E41"
$cells += "D42"
$cells += "E42"
$cells += "E43"
$cells += "E44"
$cells += "D45"
$cells += "E45"
$cells += "D46"
$cells += "E46"
$cells += "D47"
$cells += "D48"
$cells += "E48"
$cells += "E49"
$cells += "D50"
$cells += "E50"
$cells += "D51"
$cells += "E51"

foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.611.33"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.811.95"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "226.26"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  +3.60%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "37.88"
$ws.Range("E8").Value = "  +8.32%  "
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "0.0681"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "11.38"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "1.824.08"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "34.578.83"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "244.52"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").Value = "171.84"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "17.35"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "3.82"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "1.366.94"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "0.657"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("D37").Value = "1.07"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  +8.72%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").Value = "81.14"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "14.01"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("D46").Value = "0.0503"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "1.974.92"
$ws.Range("D48").Value = "5.83"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "102.80"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  -7.01%  "

foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}
